$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.430.99"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "3.896.76"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.20"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.86"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.92%  "
$ws.Range("D7").Value = "3.896.06"
$ws.Range("E7").Value = "  +1.85%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.40"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("E12").Value = "  +2.16%  "
$ws.Range("E13").Value = "  +4.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.24"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.96%  "
$ws.Range("D15").Value = "4.507.54"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").Value = "3.886.45"
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("D17").Value = "69.485.68"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.79"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +10.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.63"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.07"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "489.57"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.744"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000167"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.29"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.30"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.39"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.16"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.77%  "
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.99"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.38%  "
$ws.Range("D31").Value = "4.044.47"
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.82"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.97"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("D35").Value = "3.864.66"
$ws.Range("E35").Value = "  +2.49%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.12"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.49%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.143"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.40"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +15.17%  "
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.326"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.08"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "437.64"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.04"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.70"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.40%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000278"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +23.52%  "
$ws.Range("E49").Value = "  +2.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.34"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "141.40"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.68%  "
